$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Regular text-value updates (dates, measurements, temperatures, pressures, wind) ---
$ws.Range("E2").Value = "2026-02-16 19:18:28"
$ws.Range("I2").Value = "20.2 mm"
$ws.Range("E3").Value = "2026-02-16 19:18:30"
$ws.Range("G3").Value = "235 cm"
$ws.Range("E4").Value = "2026-02-16 19:18:33"
$ws.Range("E5").Value = "2026-02-16 19:18:35"
$ws.Range("I5").Value = "24.1 mm"
$ws.Range("N5").Value = "-1.6 °C 18:59 TU"
$ws.Range("E6").Value = "2026-02-16 19:18:38"
$ws.Range("O6").Value = "11.8 °C"
$ws.Range("E7").Value = "2026-02-16 19:18:40"
$ws.Range("J7").Value = "1013.4 hPa"
$ws.Range("O7").Value = "16.4 °C"
$ws.Range("E8").Value = "2026-02-16 19:18:43"
$ws.Range("J8").Value = "1013.0 hPa"
$ws.Range("O8").Value = "12.6 °C"
$ws.Range("E9").Value = "2026-02-16 19:18:45"
$ws.Range("O9").Value = "11.3 °C"
$ws.Range("E10").Value = "2026-02-16 19:18:48"
$ws.Range("O10").Value = "11.0 °C"
$ws.Range("E11").Value = "2026-02-16 19:18:50"
$ws.Range("E12").Value = "2026-02-16 19:18:52"
$ws.Range("O12").Value = "10.8 °C"
$ws.Range("E13").Value = "2026-02-16 19:18:55"
$ws.Range("J13").Value = "1014.8 hPa"
$ws.Range("O13").Value = "5.7 °C"
$ws.Range("E14").Value = "2026-02-16 19:18:57"
$ws.Range("E15").Value = "2026-02-16 19:18:59"
$ws.Range("O15").Value = "11.4 °C"
$ws.Range("E16").Value = "2026-02-16 19:19:02"
$ws.Range("E17").Value = "2026-02-16 19:19:04"
$ws.Range("E18").Value = "2026-02-16 19:19:07"
$ws.Range("E19").Value = "2026-02-16 19:19:09"
$ws.Range("O19").Value = "7.0 °C"
$ws.Range("E20").Value = "2026-02-16 19:19:12"
$ws.Range("E21").Value = "2026-02-16 19:19:14"
$ws.Range("J21").Value = "1014.3 hPa"
$ws.Range("E22").Value = "2026-02-16 19:19:16"
$ws.Range("E23").Value = "2026-02-16 19:19:19"
$ws.Range("I23").Value = "13.5 mm"
$ws.Range("E24").Value = "2026-02-16 19:19:21"
$ws.Range("J24").Value = "1016.6 hPa"
$ws.Range("E25").Value = "2026-02-16 19:19:24"
$ws.Range("I25").Value = "5.7 mm"
$ws.Range("O25").Value = "0.7 °C"
$ws.Range("E26").Value = "2026-02-16 19:19:26"
$ws.Range("E27").Value = "2026-02-16 19:19:28"
$ws.Range("E28").Value = "2026-02-16 19:19:30"
$ws.Range("E29").Value = "2026-02-16 19:19:32"
$ws.Range("O29").Value = "10.9 °C"
$ws.Range("E30").Value = "2026-02-16 19:19:35"
$ws.Range("J30").Value = "1012.4 hPa"
$ws.Range("E31").Value = "2026-02-16 19:19:37"
$ws.Range("J31").Value = "1011.6 hPa"
$ws.Range("O31").Value = "14.6 °C"
$ws.Range("E32").Value = "2026-02-16 19:19:40"
$ws.Range("O32").Value = "8.6 °C"
$ws.Range("E33").Value = "2026-02-16 19:19:42"
$ws.Range("J33").Value = "1013.8 hPa"
$ws.Range("E34").Value = "2026-02-16 19:19:45"
$ws.Range("L34").Value = "63.4 km/h - 49º 18:49 TU"
$ws.Range("N34").Value = "2.1 °C 18:30 TU"
$ws.Range("E35").Value = "2026-02-16 19:19:47"
$ws.Range("E36").Value = "2026-02-16 19:19:50"
$ws.Range("J36").Value = "1012.7 hPa"
$ws.Range("O36").Value = "11.7 °C"
$ws.Range("E37").Value = "2026-02-16 19:19:53"
$ws.Range("E38").Value = "2026-02-16 19:19:55"
$ws.Range("E39").Value = "2026-02-16 19:19:58"
$ws.Range("I39").Value = "3.4 mm"
$ws.Range("N39").Value = "-1.5 °C 18:59 TU"
$ws.Range("O39").Value = "0.4 °C"
$ws.Range("E40").Value = "2026-02-16 19:20:00"
$ws.Range("J40").Value = "1016.4 hPa"
$ws.Range("O40").Value = "7.0 °C"
$ws.Range("E41").Value = "2026-02-16 19:20:03"
$ws.Range("E42").Value = "2026-02-16 19:20:05"
$ws.Range("E43").Value = "2026-02-16 19:20:08"
$ws.Range("O43").Value = "8.6 °C"
$ws.Range("E44").Value = "2026-02-16 19:20:11"
$ws.Range("I44").Value = "8.3 mm"
$ws.Range("L44").Value = "56.5 km/h - 70º 18:35 TU"
$ws.Range("O44").Value = "-0.1 °C"
$ws.Range("E45").Value = "2026-02-16 19:20:14"
$ws.Range("I45").Value = "17.3 mm"
$ws.Range("E46").Value = "2026-02-16 19:20:16"
$ws.Range("O46").Value = "16.1 °C"

# --- Percentage cells: must stay as literal text "NN%" (not auto-converted to a number) ---
# Force text number format, assign the text, then restore the original border-only style
# by pasting formats from an untouched same-style cell (H2), so the cell keeps its original
# style index (s="3") instead of acquiring a new text-formatted style.
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "59%"
$ws.Range("H2").Copy() | Out-Null
$ws.Range("H4").PasteSpecial(-4122) | Out-Null

$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "80%"
$ws.Range("H2").Copy() | Out-Null
$ws.Range("H12").PasteSpecial(-4122) | Out-Null

$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "68%"
$ws.Range("H2").Copy() | Out-Null
$ws.Range("H15").PasteSpecial(-4122) | Out-Null

$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "70%"
$ws.Range("H2").Copy() | Out-Null
$ws.Range("H17").PasteSpecial(-4122) | Out-Null

$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "71%"
$ws.Range("H2").Copy() | Out-Null
$ws.Range("H21").PasteSpecial(-4122) | Out-Null

$ws.Range("H24").NumberFormat = "@"
$ws.Range("H24").Value = "71%"
$ws.Range("H2").Copy() | Out-Null
$ws.Range("H24").PasteSpecial(-4122) | Out-Null

$ws.Range("H25").NumberFormat = "@"
$ws.Range("H25").Value = "87%"
$ws.Range("H2").Copy() | Out-Null
$ws.Range("H25").PasteSpecial(-4122) | Out-Null

$ws.Range("H35").NumberFormat = "@"
$ws.Range("H35").Value = "72%"
$ws.Range("H2").Copy() | Out-Null
$ws.Range("H35").PasteSpecial(-4122) | Out-Null

$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = "74%"
$ws.Range("H2").Copy() | Out-Null
$ws.Range("H36").PasteSpecial(-4122) | Out-Null

$ws.Range("H38").NumberFormat = "@"
$ws.Range("H38").Value = "70%"
$ws.Range("H2").Copy() | Out-Null
$ws.Range("H38").PasteSpecial(-4122) | Out-Null

$ws.Range("H39").NumberFormat = "@"
$ws.Range("H39").Value = "78%"
$ws.Range("H2").Copy() | Out-Null
$ws.Range("H39").PasteSpecial(-4122) | Out-Null

$ws.Range("H46").NumberFormat = "@"
$ws.Range("H46").Value = "56%"
$ws.Range("H2").Copy() | Out-Null
$ws.Range("H46").PasteSpecial(-4122) | Out-Null

